$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "TXT" format block, added below the existing "LIN" example (rows 20-28).
$ws.Range("A20").Value = "TXT"

$ws.Range("A21").Value = "T"

$ws.Range("A22").Value = "[empty]"

$ws.Range("A23").Value = "D"
$ws.Range("B23").Value = "As ""June 2000"""

$ws.Range("A24").Value = "L"
$ws.Range("B24").Value = "With comma, not :"

$ws.Range("A25").Value = "E"

$ws.Range("A26").Value = "S"
$ws.Range("B26").Value = "As ""Semifinal, Segment 1"""

$ws.Range("A27").Value = "F"
$ws.Range("B27").Value = "not given?"

$ws.Range("A28").Value = "K"
$ws.Range("B28").Value = "Meltzer (470) vs. Schwarz (451)"

# Move the active selection to match the edited workbook's cursor state.
$ws.Range("B19").Select()
